$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: append rows 34-39
# ---------------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-01", "18:03:21", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "18:03:22", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "18:03:26", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "18:03:31", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "18:03:36", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-01", "18:03:41", "18:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 34
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $values = $pirRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $text = $values[$c]
        $cell = $pir.Cells.Item($r, $col)
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 32-36
# ---------------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-01", "18:03:21", "18:00", "Bathroom", "81.1%", "Active"),
    @("2026-02-01", "18:03:22", "18:00", "Bathroom", "80.0%", "Active"),
    @("2026-02-01", "18:03:26", "18:00", "Bathroom", "80.9%", "Active"),
    @("2026-02-01", "18:03:36", "18:00", "Bathroom", "81.0%", "Active"),
    @("2026-02-01", "18:03:41", "18:00", "Bathroom", "79.9%", "Active")
)

$startRow = 32
for ($i = 0; $i -lt $humidityRows.Count; $i++) {
    $r = $startRow + $i
    $values = $humidityRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $text = $values[$c]
        $cell = $humidity.Cells.Item($r, $col)
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 32-36
# ---------------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-01", "18:03:21", "18:00", "Bathroom", "28.8C", "Active"),
    @("2026-02-01", "18:03:22", "18:00", "Bathroom", "28.8C", "Active"),
    @("2026-02-01", "18:03:27", "18:00", "Bathroom", "28.8C", "Active"),
    @("2026-02-01", "18:03:37", "18:00", "Bathroom", "28.9C", "Active"),
    @("2026-02-01", "18:03:42", "18:00", "Bathroom", "28.9C", "Active")
)

$startRow = 32
for ($i = 0; $i -lt $temperatureRows.Count; $i++) {
    $r = $startRow + $i
    $values = $temperatureRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $text = $values[$c]
        $cell = $temperature.Cells.Item($r, $col)
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    }
}

Write-Output "Applied PIR/Humidity/Temperature updates"
